$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.26"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-4.45%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.25"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-8.40%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.079"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.46%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07676"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-6.13%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.243"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.73%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.608"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-10.71%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9149"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-3.63%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1034"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-8.28%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1744"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-6.84%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09000"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-3.94%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04428"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.21%"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.37%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001257"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-3.66%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005836"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.69%"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2,414.78%"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.05%"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-4.68%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3308"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.65%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.019"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-5.96%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1348"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.95%"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.22%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04139"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.06%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001205"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-3.62%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004100"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-4.16%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001300"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "8.31%"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02393"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-10.29%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05183"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-6.77%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007917"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.74%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1316"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-5.99%"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-10.79%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001949"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-7.60%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008370"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "9.22%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3333"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "4.04%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006420"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.76%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.12%"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-26.74%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004244"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "36.60%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.12%"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.12%"
